$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# The paragraph "ng g c -skipTests true | to prevent creation of test
# [Wingdings arrow] New version to do this" has its run boundaries
# reshuffled (adjacent same-formatted runs collapse together) while the
# visible text stays identical. Doing a no-op Find/Replace across the
# old run boundaries forces Word to re-flow/merge the runs exactly like
# the authoring edit did.
$d.Content.Find.Execute(
    " true | to prevent creation of test ", $true, $false, $false, $false,
    $false, $true, 1, $false,
    " true | to prevent creation of test ", 2) | Out-Null

# --- Change 2 --------------------------------------------------------
# The last (empty) paragraph in the document body gets a new run of
# text "--spec false", formatted with the Arial font (ascii/hAnsi/cs),
# matching the paragraph mark's run properties already on that
# paragraph.
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertAfter("--spec false")
$newLast = $d.Paragraphs($d.Paragraphs.Count)
$newLast.Range.Font.Name = "Arial"
$newLast.Range.Font.NameBi = "Arial"
